$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old table (rows 1-7, columns A-I) entirely; the new layout
# replaces it with a narrower table (A-E) plus a legend block further down.
$ws.Range("A1:I7").ClearContents()

# --- Header rows -----------------------------------------------------
# Row 1: column "type" markers (all "int")
$ws.Cells.Item(1,1).Value = "int"
$ws.Cells.Item(1,2).Value = "int"
$ws.Cells.Item(1,3).Value = "int"
$ws.Cells.Item(1,4).Value = "int"
$ws.Cells.Item(1,5).Value = "int"

# Row 2: column names. "ID" already existed in the shared strings table;
# set it first so the legend strings below (which also reuse / introduce
# shared strings) land in the order seen in the target file.
$ws.Cells.Item(2,1).Value = "ID"

# --- Legend block (rows 24-28) ---------------------------------------
# These reuse / introduce new shared strings; the order they are first
# assigned here controls their position in the shared string table.
$ws.Cells.Item(26,1).Value = "#1 : Blue"
$ws.Cells.Item(25,1).Value = "#0 : Red"
$ws.Cells.Item(27,1).Value = "#2 : Green"
$ws.Cells.Item(28,1).Value = "#3 : Black"

# Finish row 2 headers (new shared strings, appended after the legend ones)
$ws.Cells.Item(2,2).Value = "jumpPower"
$ws.Cells.Item(2,3).Value = "mass"
$ws.Cells.Item(2,4).Value = "moveSpeed"
$ws.Cells.Item(2,5).Value = "attackSpeed"

# Row 24: legend header, reuses existing shared string "# ID"
$ws.Cells.Item(24,1).Value = "# ID"

# --- Data rows 3-6 -----------------------------------------------------
$ws.Cells.Item(3,1).Value = 0
$ws.Cells.Item(3,2).Value = 10
$ws.Cells.Item(3,3).Value = 7
$ws.Cells.Item(3,4).Value = 8
$ws.Cells.Item(3,5).Value = 5

$ws.Cells.Item(4,1).Value = 1
$ws.Cells.Item(4,2).Value = 10
$ws.Cells.Item(4,3).Value = 5
$ws.Cells.Item(4,4).Value = 10
$ws.Cells.Item(4,5).Value = 7

$ws.Cells.Item(5,1).Value = 2
$ws.Cells.Item(5,2).Value = 10
$ws.Cells.Item(5,3).Value = 10
$ws.Cells.Item(5,4).Value = 5
$ws.Cells.Item(5,5).Value = 3

$ws.Cells.Item(6,1).Value = 3
$ws.Cells.Item(6,2).Value = 10
$ws.Cells.Item(6,3).Value = 6
$ws.Cells.Item(6,4).Value = 7
$ws.Cells.Item(6,5).Value = 6

$ws.Range("J14").Select()
